$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "22-10-2023"
$ws.Range("B14").Value = "Docker"
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("C14").Value = "        Dot_Net_Application"
$ws.Range("D14").Value = "1.Download Docker"
$ws.Range("D15").Value = "2.wsl2 required     wsl --install"
$ws.Range("D16").Value = "3.enable Virulazation      bypressing f12"
$ws.Range("D17").Value = "4.for checking taskmanager performance"
$ws.Range("D18").Value = "5.window 21h1 not less than this"
